$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.185.91'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.851.71'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.37'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4702'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.48%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06550'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07964'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '97.54'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.855.06'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.093'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6752'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '269.07'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -3.70%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.149.88'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.62'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +7.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007650'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +4.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.094.26'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.000'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.201'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -5.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.135'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '167.18'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.155'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.82'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.929'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09848'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.463'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.281'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.992'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -2.67%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04695'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.116'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6982'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.18%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01868'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +2.86%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.318'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '73.13'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.930'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9992'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8385'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '103.16'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4130'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '936.19'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.132'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.018'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -2.18%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.85'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.85%  '
